$wb = $excel.ActiveWorkbook

# --- 1. "Snow clearing" sheet: move the saved selection from F4 to G47 ---
#        (must happen before the new sheet is added/activated below, since
#         adding a sheet makes it the active tab)
$wsSnowClearing = $wb.Worksheets.Item("Snow clearing")
[void]$wsSnowClearing.Activate()
$wsSnowClearing.Range("G47").Select() | Out-Null

# --- 2. Add the new "ESTM Coefficients" sheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ESTM Coefficients"

$newSheet.Range("A1").Value = "SUEWS_ESTMCoefficients"
$newSheet.Range("B1").Value = "SUEWS_ESTMCoefficients.txt"
$newSheet.Range("C1").Value = "ESTM Coefficients (codes will be added to the final output according to file input not what is selected in this tab)"
$newSheet.Range("F1").Value = 81

$newSheet.Range("F7").Select() | Out-Null

# --- 3. "Water Use (Automatic)" sheet: give columns A & B explicit widths ---
$wsWaterAuto = $wb.Worksheets.Item("Water Use (Automatic)")
$wsWaterAuto.Columns.Item(1).ColumnWidth = 20 + 2/3
$wsWaterAuto.Columns.Item(2).ColumnWidth = 29 + 5/6
